# The workbook's daily price log gained one new record: a row is inserted
# right before the existing row 30, pushing every subsequent record (old
# rows 30-126) down by one (to new rows 31-127). The new row 30 carries a
# fresh observation for "Arándano (blue)" at the Macroferia Regional de
# Talca market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 30; Excel shifts rows 30:126 down to
# 31:127 and extends the used range/dimension accordingly.
$ws.Rows(30).Insert()

# Populate the newly inserted row 30 with the new record's data.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44972
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100101
$ws.Range("H30").Value = "Berries"
$ws.Range("I30").Value = 100101001
$ws.Range("J30").Value = "Arándano (blue)"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 3000
$ws.Range("O30").Value = 3000
$ws.Range("P30").Value = 3000
$ws.Range("Q30").Value = "$/bandeja 2 kilos"
$ws.Range("R30").Value = "Provincia de Curicó"
$ws.Range("S30").Value = 1500
$ws.Range("T30").Value = 2
